# Add team record (Wins / Losses / Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the existing data before we add new columns.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count()
$lastDataCol = $used.Columns.Count()   # AC = 29 in this workbook

$headerRow = 1
$newCol1 = $lastDataCol + 1   # AD - Wins
$newCol2 = $lastDataCol + 2   # AE - Losses
$newCol3 = $lastDataCol + 3   # AF - Ties

# --- Header row: new columns AD, AE, AF ---
# Copy the formatting of the existing last header cell (bold, centered,
# thin border) onto the new header cells before setting their text.
$lastHeaderCell = $ws.Cells.Item($headerRow, $lastDataCol)
$lastHeaderCell.Copy() | Out-Null
$newHeaderRange = $ws.Range($ws.Cells.Item($headerRow, $newCol1), $ws.Cells.Item($headerRow, $newCol3))
$newHeaderRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($headerRow, $newCol1).Value = "Wins"
$ws.Cells.Item($headerRow, $newCol2).Value = "Losses"
$ws.Cells.Item($headerRow, $newCol3).Value = "Ties"

# --- Data rows: team record values for every player ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newCol1).Value = 71
    $ws.Cells.Item($r, $newCol2).Value = 91
    $ws.Cells.Item($r, $newCol3).Value = 0
}
